$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2;  D=45001; M=60; N=17000; O=18000; P=17500; Q="`$/caja 18 kilos"; R="Región Metropolitana";   S=972;  T=18 }
    @{ Row=3;  D=45002; M=30; N=18000; O=18000; P=18000; Q="`$/caja 18 kilos"; R="Región Metropolitana";   S=1000; T=18 }
    @{ Row=4;  D=45050; M=40; N=14000; O=14000; P=14000; Q="`$/caja 18 kilos"; R="Región Metropolitana";   S=778;  T=18 }
    @{ Row=5;  D=45037; M=60; N=16000; O=16000; P=16000; Q="`$/caja 18 kilos"; R="Región Metropolitana";   S=889;  T=18 }
    @{ Row=6;  D=45099; M=40; N=22000; O=22000; P=22000; Q="`$/caja 18 kilos"; R="Región Metropolitana";   S=1222; T=18 }
    @{ Row=7;  D=45036; M=60; N=15000; O=16000; P=15500; Q="`$/caja 18 kilos"; R="Región Metropolitana";   S=861;  T=18 }
    @{ Row=8;  D=45028; M=50; N=18000; O=18000; P=18000; Q="`$/caja 18 kilos"; R="Región Metropolitana";   S=1000; T=18 }
    @{ Row=9;  D=45049; M=80; N=15000; O=15000; P=15000; Q="`$/caja 18 kilos"; R="Región Metropolitana";   S=833;  T=18 }
    @{ Row=10; D=45021; M=60; N=15000; O=16000; P=15500; Q="`$/caja 18 kilos"; R="Provincia de Los Andes"; S=861;  T=18 }
    @{ Row=11; D=45041; M=60; N=15000; O=15000; P=15000; Q="`$/caja 18 kilos"; R="Región Metropolitana";   S=833;  T=18 }
    @{ Row=12; D=45030; M=40; N=18000; O=18000; P=18000; Q="`$/caja 18 kilos"; R="Región Metropolitana";   S=1000; T=18 }
    @{ Row=13; D=45096; M=50; N=23000; O=23000; P=23000; Q="`$/caja 18 kilos"; R="Región Metropolitana";   S=1278; T=18 }
    @{ Row=14; D=45062; M=90; N=13000; O=14000; P=13444; Q="`$/caja 18 kilos"; R="Región Metropolitana";   S=747;  T=18 }
    @{ Row=15; D=45014; M=30; N=18000; O=18000; P=18000; Q="`$/caja 18 kilos"; R="Región Metropolitana";   S=1000; T=18 }
    @{ Row=16; D=45043; M=60; N=15000; O=15000; P=15000; Q="`$/caja 18 kilos"; R="Región Metropolitana";   S=833;  T=18 }
    @{ Row=17; D=44999; M=60; N=17000; O=18000; P=17500; Q="`$/caja 18 kilos"; R="Región Metropolitana";   S=972;  T=18 }
    @{ Row=18; D=45020; M=50; N=15000; O=15000; P=15000; Q="`$/caja 16 kilos"; R="Provincia de Los Andes"; S=938;  T=16 }
    @{ Row=19; D=45033; M=60; N=15000; O=16000; P=15500; Q="`$/caja 18 kilos"; R="Región Metropolitana";   S=861;  T=18 }
    @{ Row=20; D=45089; M=60; N=22000; O=23000; P=22500; Q="`$/caja 18 kilos"; R="Región Metropolitana";   S=1250; T=18 }
    @{ Row=21; D=45044; M=60; N=15000; O=15000; P=15000; Q="`$/caja 18 kilos"; R="Región Metropolitana";   S=833;  T=18 }
    @{ Row=22; D=45091; M=50; N=22000; O=22000; P=22000; Q="`$/caja 18 kilos"; R="Región Metropolitana";   S=1222; T=18 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 4).Value  = $r.D    # D - Fecha
    $ws.Cells.Item($row, 13).Value = $r.M    # M - Volumen
    $ws.Cells.Item($row, 14).Value = $r.N    # N - Precio minimo
    $ws.Cells.Item($row, 15).Value = $r.O    # O - Precio maximo
    $ws.Cells.Item($row, 16).Value = $r.P    # P - Precio promedio ponderado
    $ws.Cells.Item($row, 17).Value = $r.Q    # Q - Unidad de comercializacion
    $ws.Cells.Item($row, 18).Value = $r.R    # R - Origen
    $ws.Cells.Item($row, 19).Value = $r.S    # S - Precio $/Kg
    $ws.Cells.Item($row, 20).Value = $r.T    # T - Kg / unidad
}
